$d = $word.ActiveDocument

function Replace-OneText($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 1)
}

# --- Guidance copy (first occurrence of each) ---
Replace-OneText "John Vincent" "asdf"
Replace-OneText "12-ambot" "adsf"
Replace-OneText "example teacher" "asdf"
Replace-OneText "2024-03-27   " "2024-03-20   "
Replace-OneText "  09:01:00   " "  19:08:00   "
Replace-OneText "19:57:00" "08:09:00"
Replace-OneText "                                 Example Counselor________________ " "                                 asdfasdf________________ "
Replace-OneText "                          example teacher_________ " "                          asdf_________ "

# --- Teacher's copy (second occurrence of each) ---
Replace-OneText "John Vincent " "asdf "
Replace-OneText "12-ambot" "adsf"
Replace-OneText "example teacher" "asdf"
Replace-OneText "          2024-03-27   " "          2024-03-20   "
Replace-OneText "  09:01:00   " "  19:08:00   "
Replace-OneText "19:57:00" "08:09:00"
Replace-OneText "                                 Example Counselor________________ " "                                 asdfasdf________________ "
Replace-OneText "                          example teacher_________ " "                          asdf_________ "
